# "Generate Report for Handoff"
#
# A new handback/handoff report run updates the timestamps recorded for the
# c1abf968-06d8-48fe-8350-6b42d6fadf4d file (row 7 on every sheet):
#   - Overview sheet: "Latest Handoff Date" column (D)
#   - zh-cn sheet:     "Latest Handback DateTime" column (E)
#   - de-de sheet:     "Latest Handback DateTime" column (E)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-26-14 02:26:25"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-14 02:26:22"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-14 02:26:25"
